$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 187, pushing every existing row
# (187..259) down by one (188..260), matching the committed diff.
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R repeat the same template values as
# the surrounding "Betarraga" rows; D,J,K,L,M,P carry the new figures.
$ws.Range("A187").Value = 7
$ws.Range("B187").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C187").Value = "Ñuble"
$ws.Range("D187").Value = 44468
$ws.Range("E187").Value = 16
$ws.Range("F187").Value = 100114014
$ws.Range("G187").Value = "Betarraga"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 160
$ws.Range("K187").Value = 700
$ws.Range("L187").Value = 800
$ws.Range("M187").Value = 750
$ws.Range("N187").Value = "$/paquete 5 unidades"
$ws.Range("O187").Value = "Provincia de Diguillín"
$ws.Range("P187").Value = 150
$ws.Range("Q187").Value = 5
$ws.Range("R187").Value = "Hortaliza"
